# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (Doctor Davis, "Región de O'Higgins") right before
# the existing "Early Majestic" / San Felipe de Aconcagua block that starts
# at row 369, pushing all following rows down by 3 (dimension A1:T384 -> A1:T387).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 369:371 - everything currently there (and below)
# shifts down to 372:387.
$ws.Range("A369:T371").EntireRow.Insert()

# New row 369: Doctor Davis / Especial
$ws.Cells.Item(369, 1).Value  = 8
$ws.Cells.Item(369, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(369, 3).Value  = "Coquimbo"
$ws.Cells.Item(369, 4).Value  = "2/23/2022"
$ws.Cells.Item(369, 5).Value  = 4
$ws.Cells.Item(369, 6).Value  = "Fruta"
$ws.Cells.Item(369, 7).Value  = 100103
$ws.Cells.Item(369, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(369, 9).Value  = 100103004
$ws.Cells.Item(369, 10).Value = "Durazno"
$ws.Cells.Item(369, 11).Value = "Doctor Davis"
$ws.Cells.Item(369, 12).Value = "Especial"
$ws.Cells.Item(369, 13).Value = 20
$ws.Cells.Item(369, 14).Value = 350000
$ws.Cells.Item(369, 15).Value = 360000
$ws.Cells.Item(369, 16).Value = 355000
$ws.Cells.Item(369, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(369, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(369, 19).Value = 888
$ws.Cells.Item(369, 20).Value = 400

# New row 370: Doctor Davis / Primera
$ws.Cells.Item(370, 1).Value  = 8
$ws.Cells.Item(370, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(370, 3).Value  = "Coquimbo"
$ws.Cells.Item(370, 4).Value  = "2/23/2022"
$ws.Cells.Item(370, 5).Value  = 4
$ws.Cells.Item(370, 6).Value  = "Fruta"
$ws.Cells.Item(370, 7).Value  = 100103
$ws.Cells.Item(370, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(370, 9).Value  = 100103004
$ws.Cells.Item(370, 10).Value = "Durazno"
$ws.Cells.Item(370, 11).Value = "Doctor Davis"
$ws.Cells.Item(370, 12).Value = "Primera"
$ws.Cells.Item(370, 13).Value = 16
$ws.Cells.Item(370, 14).Value = 300000
$ws.Cells.Item(370, 15).Value = 310000
$ws.Cells.Item(370, 16).Value = 305000
$ws.Cells.Item(370, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(370, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(370, 19).Value = 762
$ws.Cells.Item(370, 20).Value = 400

# New row 371: Doctor Davis / Segunda
$ws.Cells.Item(371, 1).Value  = 8
$ws.Cells.Item(371, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(371, 3).Value  = "Coquimbo"
$ws.Cells.Item(371, 4).Value  = "2/23/2022"
$ws.Cells.Item(371, 5).Value  = 4
$ws.Cells.Item(371, 6).Value  = "Fruta"
$ws.Cells.Item(371, 7).Value  = 100103
$ws.Cells.Item(371, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(371, 9).Value  = 100103004
$ws.Cells.Item(371, 10).Value = "Durazno"
$ws.Cells.Item(371, 11).Value = "Doctor Davis"
$ws.Cells.Item(371, 12).Value = "Segunda"
$ws.Cells.Item(371, 13).Value = 16
$ws.Cells.Item(371, 14).Value = 270000
$ws.Cells.Item(371, 15).Value = 280000
$ws.Cells.Item(371, 16).Value = 275000
$ws.Cells.Item(371, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(371, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(371, 19).Value = 688
$ws.Cells.Item(371, 20).Value = 400
